$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 164; all existing rows 164:241 shift down to 166:243.
$ws.Range("A164:A165").EntireRow.Insert()

# New row 164 - "Primera" quality record for date 44489 (2021-10-20)
$ws.Range("A164").Value = 8
$ws.Range("B164").Value = "Terminal La Palmera de La Serena"
$ws.Range("C164").Value = "Coquimbo"
$ws.Range("D164").Value = 44489
$ws.Range("E164").Value = 4
$ws.Range("F164").Value = 100112017
$ws.Range("G164").Value = "Apio"
$ws.Range("H164").Value = "Americana (o)"
$ws.Range("I164").Value = "Primera"
$ws.Range("J164").Value = 2000
$ws.Range("K164").Value = 6800
$ws.Range("L164").Value = 7000
$ws.Range("M164").Value = 6900
$ws.Range("N164").Value = "$/docena de matas"
$ws.Range("O164").Value = "Provincia del Elquí"
$ws.Range("P164").Value = 1150
$ws.Range("Q164").Value = 6
$ws.Range("R164").Value = "Hortaliza"

# New row 165 - "Segunda" quality record for date 44489 (2021-10-20)
$ws.Range("A165").Value = 8
$ws.Range("B165").Value = "Terminal La Palmera de La Serena"
$ws.Range("C165").Value = "Coquimbo"
$ws.Range("D165").Value = 44489
$ws.Range("E165").Value = 4
$ws.Range("F165").Value = 100112017
$ws.Range("G165").Value = "Apio"
$ws.Range("H165").Value = "Americana (o)"
$ws.Range("I165").Value = "Segunda"
$ws.Range("J165").Value = 1400
$ws.Range("K165").Value = 5500
$ws.Range("L165").Value = 6000
$ws.Range("M165").Value = 5750
$ws.Range("N165").Value = "$/docena de matas"
$ws.Range("O165").Value = "Provincia del Elquí"
$ws.Range("P165").Value = 958
$ws.Range("Q165").Value = 6
$ws.Range("R165").Value = "Hortaliza"

# Ensure the date cells retain the same date style as the rest of column D (numFmtId 165, style index 2)
$ws.Range("D164").NumberFormat = $ws.Range("D166").NumberFormat
$ws.Range("D165").NumberFormat = $ws.Range("D166").NumberFormat
